$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").PrefixCharacter = "'"
$ws.Range("D2").Value = '67.679.88'
$ws.Range("E2").PrefixCharacter = "'"
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").PrefixCharacter = "'"
$ws.Range("D3").Value = '2.638.74'
$ws.Range("E3").PrefixCharacter = "'"
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("D4").PrefixCharacter = "'"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").PrefixCharacter = "'"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").PrefixCharacter = "'"
$ws.Range("D5").Value = '595.77'
$ws.Range("E5").PrefixCharacter = "'"
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").PrefixCharacter = "'"
$ws.Range("D6").Value = '168.73'
$ws.Range("E6").PrefixCharacter = "'"
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").PrefixCharacter = "'"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").PrefixCharacter = "'"
$ws.Range("E8").Value = '  -2.05%  '
$ws.Range("D9").PrefixCharacter = "'"
$ws.Range("D9").Value = '2.637.56'
$ws.Range("E9").PrefixCharacter = "'"
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("D10").PrefixCharacter = "'"
$ws.Range("D10").Value = '0.141'
$ws.Range("E10").PrefixCharacter = "'"
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("E11").PrefixCharacter = "'"
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("E12").PrefixCharacter = "'"
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("E13").PrefixCharacter = "'"
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").PrefixCharacter = "'"
$ws.Range("D14").Value = '27.72'
$ws.Range("E14").PrefixCharacter = "'"
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").PrefixCharacter = "'"
$ws.Range("D15").Value = '3.111.75'
$ws.Range("E15").PrefixCharacter = "'"
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").PrefixCharacter = "'"
$ws.Range("D16").Value = '0.0000183'
$ws.Range("E16").PrefixCharacter = "'"
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").PrefixCharacter = "'"
$ws.Range("D17").Value = '67.652.04'
$ws.Range("E17").PrefixCharacter = "'"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").PrefixCharacter = "'"
$ws.Range("D18").Value = '2.635.79'
$ws.Range("E18").PrefixCharacter = "'"
$ws.Range("E18").Value = '  -1.71%  '
$ws.Range("D19").PrefixCharacter = "'"
$ws.Range("D19").Value = '12.10'
$ws.Range("E19").PrefixCharacter = "'"
$ws.Range("E19").Value = '  +3.04%  '
$ws.Range("D20").PrefixCharacter = "'"
$ws.Range("D20").Value = '8.08'
$ws.Range("E20").PrefixCharacter = "'"
$ws.Range("E20").Value = '  +2.49%  '
$ws.Range("D21").PrefixCharacter = "'"
$ws.Range("D21").Value = '358.75'
$ws.Range("E21").PrefixCharacter = "'"
$ws.Range("E21").Value = '  -1.77%  '
$ws.Range("E22").PrefixCharacter = "'"
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").PrefixCharacter = "'"
$ws.Range("D23").Value = '4.72'
$ws.Range("E23").PrefixCharacter = "'"
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("D24").PrefixCharacter = "'"
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").PrefixCharacter = "'"
$ws.Range("E24").Value = '  -3.93%  '
$ws.Range("D25").PrefixCharacter = "'"
$ws.Range("D25").Value = '10.40'
$ws.Range("E25").PrefixCharacter = "'"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").PrefixCharacter = "'"
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").PrefixCharacter = "'"
$ws.Range("D27").Value = '69.99'
$ws.Range("E27").PrefixCharacter = "'"
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("D28").PrefixCharacter = "'"
$ws.Range("D28").Value = '2.770.17'
$ws.Range("E28").PrefixCharacter = "'"
$ws.Range("E28").Value = '  -2.27%  '
$ws.Range("E29").PrefixCharacter = "'"
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").PrefixCharacter = "'"
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("D31").PrefixCharacter = "'"
$ws.Range("D31").Value = '550.70'
$ws.Range("E31").PrefixCharacter = "'"
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").PrefixCharacter = "'"
$ws.Range("D32").Value = '7.99'
$ws.Range("E32").PrefixCharacter = "'"
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").PrefixCharacter = "'"
$ws.Range("D33").Value = '1.37'
$ws.Range("E33").PrefixCharacter = "'"
$ws.Range("E33").Value = '  -2.37%  '
$ws.Range("D34").PrefixCharacter = "'"
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").PrefixCharacter = "'"
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("E35").PrefixCharacter = "'"
$ws.Range("E35").Value = '  +4.96%  '
$ws.Range("D36").PrefixCharacter = "'"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").PrefixCharacter = "'"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").PrefixCharacter = "'"
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("D38").PrefixCharacter = "'"
$ws.Range("D38").Value = '157.75'
$ws.Range("E38").PrefixCharacter = "'"
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("D39").PrefixCharacter = "'"
$ws.Range("D39").Value = '19.08'
$ws.Range("E39").PrefixCharacter = "'"
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("E40").PrefixCharacter = "'"
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("B41").PrefixCharacter = "'"
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").PrefixCharacter = "'"
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").PrefixCharacter = "'"
$ws.Range("D41").Value = '1.82'
$ws.Range("E41").PrefixCharacter = "'"
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("B42").PrefixCharacter = "'"
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").PrefixCharacter = "'"
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").PrefixCharacter = "'"
$ws.Range("D42").Value = '5.25'
$ws.Range("E42").PrefixCharacter = "'"
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("B43").PrefixCharacter = "'"
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").PrefixCharacter = "'"
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").PrefixCharacter = "'"
$ws.Range("D43").Value = '18.31'
$ws.Range("E43").PrefixCharacter = "'"
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("E44").PrefixCharacter = "'"
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").PrefixCharacter = "'"
$ws.Range("D45").Value = '2.45'
$ws.Range("E45").PrefixCharacter = "'"
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("E46").PrefixCharacter = "'"
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").PrefixCharacter = "'"
$ws.Range("D47").Value = '153.45'
$ws.Range("E47").PrefixCharacter = "'"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").PrefixCharacter = "'"
$ws.Range("D48").Value = '0.584'
$ws.Range("E48").PrefixCharacter = "'"
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("E49").PrefixCharacter = "'"
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("E50").PrefixCharacter = "'"
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("E51").PrefixCharacter = "'"
$ws.Range("E51").Value = '  -1.22%  '
